$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12, pushing existing rows (Explicit Text: and below) down by one.
$ws.Rows("12").Insert()

# Populate the newly inserted row with the "Large Double Number:" test case.
$ws.Range("B12").Value = "Large Double Number:"
$ws.Range("C12").Value = [double]"9.999E+307"
